$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Cells.Item(2,9).Value = [double]"0.1101680340964018"
$ws.Cells.Item(2,10).Value = [double]"0.1101680340964017"
$ws.Cells.Item(2,13).Value = [double]"0.2087793333333333"
$ws.Cells.Item(2,14).Value = [double]"0.626338"
$ws.Cells.Item(2,15).Value = [double]"0.02275344108115409"
$ws.Cells.Item(2,16).Value = [double]"0.02275344108115409"
$ws.Cells.Item(2,17).Value = [double]"0.03200935145555556"
$ws.Cells.Item(2,18).Value = [double]"0.2880841631"
$ws.Cells.Item(2,19).Value = [double]"0.002506701872839052"
$ws.Cells.Item(2,20).Value = [double]"0.002506701872839052"

# Row 3 updates
$ws.Cells.Item(3,9).Value = [double]"0.1101680340964018"
$ws.Cells.Item(3,10).Value = [double]"0.1101680340964017"
$ws.Cells.Item(3,15).Value = [double]"0.9204452022087118"
$ws.Cells.Item(3,16).Value = [double]"0.920445202208712"
$ws.Cells.Item(3,19).Value = [double]"0.1014036384207988"
$ws.Cells.Item(3,20).Value = [double]"0.1014036384207988"

# Row 4 updates
$ws.Cells.Item(4,9).Value = [double]"0.1101680340964018"
$ws.Cells.Item(4,10).Value = [double]"0.1101680340964017"
$ws.Cells.Item(4,13).Value = [double]"0.05870933333333334"
$ws.Cells.Item(4,14).Value = [double]"0.176128"
$ws.Cells.Item(4,15).Value = [double]"0.006398331365399365"
$ws.Cells.Item(4,16).Value = [double]"0.006398331365399365"
$ws.Cells.Item(4,17).Value = [double]"0.00900111928888889"
$ws.Cells.Item(4,18).Value = [double]"0.08101007360000001"
$ws.Cells.Item(4,19).Value = [double]"0.0007048915880233941"
$ws.Cells.Item(4,20).Value = [double]"0.000704891588023394"

# Row 5 updates
$ws.Cells.Item(5,9).Value = [double]"0.1101680340964018"
$ws.Cells.Item(5,10).Value = [double]"0.1101680340964017"
$ws.Cells.Item(5,13).Value = [double]"0.4594193333333333"
$ws.Cells.Item(5,14).Value = [double]"1.378258"
$ws.Cells.Item(5,15).Value = [double]"0.05006899181852175"
$ws.Cells.Item(5,16).Value = [double]"0.05006899181852175"
$ws.Cells.Item(5,17).Value = [double]"0.0704366407888889"
$ws.Cells.Item(5,18).Value = [double]"0.6339297671"
$ws.Cells.Item(5,19).Value = [double]"0.005516002397835365"
$ws.Cells.Item(5,20).Value = [double]"0.005516002397835364"

# Row 6 updates
$ws.Cells.Item(6,9).Value = [double]"0.1101680340964018"
$ws.Cells.Item(6,10).Value = [double]"0.1101680340964017"
$ws.Cells.Item(6,13).Value = [double]"0.003065"
$ws.Cells.Item(6,14).Value = [double]"0.009195"
$ws.Cells.Item(6,15).Value = [double]"0.0003340335262130221"
$ws.Cells.Item(6,16).Value = [double]"0.0003340335262130222"
$ws.Cells.Item(6,17).Value = [double]"0.0004699155833333334"
$ws.Cells.Item(6,18).Value = [double]"0.00422924025"
$ws.Cells.Item(6,19).Value = [double]"3.679981690517753e-05"
$ws.Cells.Item(6,20).Value = [double]"3.679981690517753e-05"

# Row 7 updates
$ws.Cells.Item(7,9).Value = [double]"0.6869455370628789"
$ws.Cells.Item(7,10).Value = [double]"0.6869455370628788"
$ws.Cells.Item(7,13).Value = [double]"0.2087793333333333"
$ws.Cells.Item(7,14).Value = [double]"0.626338"
$ws.Cells.Item(7,15).Value = [double]"0.02275344108115409"
$ws.Cells.Item(7,16).Value = [double]"0.02275344108115409"
$ws.Cells.Item(7,17).Value = [double]"0.1995922075493334"
$ws.Cells.Item(7,18).Value = [double]"1.796329867944"
$ws.Cells.Item(7,19).Value = [double]"0.01563037480352197"
$ws.Cells.Item(7,20).Value = [double]"0.01563037480352196"

# Row 8 updates
$ws.Cells.Item(8,9).Value = [double]"0.6869455370628789"
$ws.Cells.Item(8,10).Value = [double]"0.6869455370628788"
$ws.Cells.Item(8,15).Value = [double]"0.9204452022087118"
$ws.Cells.Item(8,16).Value = [double]"0.920445202208712"
$ws.Cells.Item(8,19).Value = [double]"0.6322957237682137"
$ws.Cells.Item(8,20).Value = [double]"0.6322957237682136"

# Row 9 updates
$ws.Cells.Item(9,9).Value = [double]"0.6869455370628789"
$ws.Cells.Item(9,10).Value = [double]"0.6869455370628788"
$ws.Cells.Item(9,13).Value = [double]"0.05870933333333334"
$ws.Cells.Item(9,14).Value = [double]"0.176128"
$ws.Cells.Item(9,15).Value = [double]"0.006398331365399365"
$ws.Cells.Item(9,16).Value = [double]"0.006398331365399365"
$ws.Cells.Item(9,17).Value = [double]"0.05612588782933334"
$ws.Cells.Item(9,18).Value = [double]"0.5051329904640001"
$ws.Cells.Item(9,19).Value = [double]"0.00439530517611053"
$ws.Cells.Item(9,20).Value = [double]"0.004395305176110529"

# Row 10 updates
$ws.Cells.Item(10,9).Value = [double]"0.6869455370628789"
$ws.Cells.Item(10,10).Value = [double]"0.6869455370628788"
$ws.Cells.Item(10,13).Value = [double]"0.4594193333333333"
$ws.Cells.Item(10,14).Value = [double]"1.378258"
$ws.Cells.Item(10,15).Value = [double]"0.05006899181852175"
$ws.Cells.Item(10,16).Value = [double]"0.05006899181852175"
$ws.Cells.Item(10,17).Value = [double]"0.4392030449893334"
$ws.Cells.Item(10,18).Value = [double]"3.952827404904001"
$ws.Cells.Item(10,19).Value = [double]"0.03439467047497131"
$ws.Cells.Item(10,20).Value = [double]"0.03439467047497131"

# Row 11 updates
$ws.Cells.Item(11,9).Value = [double]"0.6869455370628789"
$ws.Cells.Item(11,10).Value = [double]"0.6869455370628788"
$ws.Cells.Item(11,13).Value = [double]"0.003065"
$ws.Cells.Item(11,14).Value = [double]"0.009195"
$ws.Cells.Item(11,15).Value = [double]"0.0003340335262130221"
$ws.Cells.Item(11,16).Value = [double]"0.0003340335262130222"
$ws.Cells.Item(11,17).Value = [double]"0.002930127740000001"
$ws.Cells.Item(11,18).Value = [double]"0.02637114966"
$ws.Cells.Item(11,19).Value = [double]"0.0002294628400614117"
$ws.Cells.Item(11,20).Value = [double]"0.0002294628400614117"

# Row 12 updates
$ws.Cells.Item(12,7).Value = [double]"0.1265133333333333"
$ws.Cells.Item(12,8).Value = [double]"0.37954"
$ws.Cells.Item(12,9).Value = [double]"0.09090808927263468"
$ws.Cells.Item(12,10).Value = [double]"0.09090808927263468"
$ws.Cells.Item(12,13).Value = [double]"0.2087793333333333"
$ws.Cells.Item(12,14).Value = [double]"0.626338"
$ws.Cells.Item(12,15).Value = [double]"0.02275344108115409"
$ws.Cells.Item(12,16).Value = [double]"0.02275344108115409"
$ws.Cells.Item(12,17).Value = [double]"0.02641336939111111"
$ws.Cells.Item(12,18).Value = [double]"0.23772032452"
$ws.Cells.Item(12,19).Value = [double]"0.002068471853065189"
$ws.Cells.Item(12,20).Value = [double]"0.002068471853065189"

# Row 13 updates
$ws.Cells.Item(13,7).Value = [double]"0.1265133333333333"
$ws.Cells.Item(13,8).Value = [double]"0.37954"
$ws.Cells.Item(13,9).Value = [double]"0.09090808927263468"
$ws.Cells.Item(13,10).Value = [double]"0.09090808927263468"
$ws.Cells.Item(13,15).Value = [double]"0.9204452022087118"
$ws.Cells.Item(13,16).Value = [double]"0.920445202208712"
$ws.Cells.Item(13,17).Value = [double]"1.068500322368889"
$ws.Cells.Item(13,18).Value = [double]"9.61650290132"
$ws.Cells.Item(13,19).Value = [double]"0.08367591461295786"
$ws.Cells.Item(13,20).Value = [double]"0.08367591461295787"

# Row 14 updates
$ws.Cells.Item(14,7).Value = [double]"0.1265133333333333"
$ws.Cells.Item(14,8).Value = [double]"0.37954"
$ws.Cells.Item(14,9).Value = [double]"0.09090808927263468"
$ws.Cells.Item(14,10).Value = [double]"0.09090808927263468"
$ws.Cells.Item(14,13).Value = [double]"0.05870933333333334"
$ws.Cells.Item(14,14).Value = [double]"0.176128"
$ws.Cells.Item(14,15).Value = [double]"0.006398331365399365"
$ws.Cells.Item(14,16).Value = [double]"0.006398331365399365"
$ws.Cells.Item(14,17).Value = [double]"0.007427513457777779"
$ws.Cells.Item(14,18).Value = [double]"0.06684762112"
$ws.Cells.Item(14,19).Value = [double]"0.000581660078961624"
$ws.Cells.Item(14,20).Value = [double]"0.000581660078961624"

# Row 15 updates
$ws.Cells.Item(15,7).Value = [double]"0.1265133333333333"
$ws.Cells.Item(15,8).Value = [double]"0.37954"
$ws.Cells.Item(15,9).Value = [double]"0.09090808927263468"
$ws.Cells.Item(15,10).Value = [double]"0.09090808927263468"
$ws.Cells.Item(15,13).Value = [double]"0.4594193333333333"
$ws.Cells.Item(15,14).Value = [double]"1.378258"
$ws.Cells.Item(15,15).Value = [double]"0.05006899181852175"
$ws.Cells.Item(15,16).Value = [double]"0.05006899181852175"
$ws.Cells.Item(15,17).Value = [double]"0.05812267125777778"
$ws.Cells.Item(15,18).Value = [double]"0.52310404132"
$ws.Cells.Item(15,19).Value = [double]"0.00455167637802899"
$ws.Cells.Item(15,20).Value = [double]"0.00455167637802899"

# Row 16 updates
$ws.Cells.Item(16,7).Value = [double]"0.1265133333333333"
$ws.Cells.Item(16,8).Value = [double]"0.37954"
$ws.Cells.Item(16,9).Value = [double]"0.09090808927263468"
$ws.Cells.Item(16,10).Value = [double]"0.09090808927263468"
$ws.Cells.Item(16,13).Value = [double]"0.003065"
$ws.Cells.Item(16,14).Value = [double]"0.009195"
$ws.Cells.Item(16,15).Value = [double]"0.0003340335262130221"
$ws.Cells.Item(16,16).Value = [double]"0.0003340335262130222"
$ws.Cells.Item(16,17).Value = [double]"0.0003877633666666667"
$ws.Cells.Item(16,18).Value = [double]"0.0034898703"
$ws.Cells.Item(16,19).Value = [double]"3.036634962102637e-05"
$ws.Cells.Item(16,20).Value = [double]"3.036634962102637e-05"

# Row 17 (new)
$ws.Cells.Item(17,1).Value = "Resolving-Mac"
$ws.Cells.Item(17,2).Value = "Efnb3"
$ws.Cells.Item(17,3).Value = "Ephb2"
$ws.Cells.Item(17,4).Value = "ECs"
$ws.Cells.Item(17,5).Value = 1
$ws.Cells.Item(17,6).Value = [double]"0.3333333333333333"
$ws.Cells.Item(17,7).Value = [double]"0.155836"
$ws.Cells.Item(17,8).Value = [double]"0.467508"
$ws.Cells.Item(17,9).Value = [double]"0.1119783395680848"
$ws.Cells.Item(17,10).Value = [double]"0.1119783395680848"
$ws.Cells.Item(17,11).Value = 2
$ws.Cells.Item(17,12).Value = [double]"0.6666666666666666"
$ws.Cells.Item(17,13).Value = [double]"0.2087793333333333"
$ws.Cells.Item(17,14).Value = [double]"0.626338"
$ws.Cells.Item(17,15).Value = [double]"0.02275344108115409"
$ws.Cells.Item(17,16).Value = [double]"0.02275344108115409"
$ws.Cells.Item(17,17).Value = [double]"0.03253533618933333"
$ws.Cells.Item(17,18).Value = [double]"0.292818025704"
$ws.Cells.Item(17,19).Value = [double]"0.002547892551727883"
$ws.Cells.Item(17,20).Value = [double]"0.002547892551727882"

# Row 18 (new)
$ws.Cells.Item(18,1).Value = "Resolving-Mac"
$ws.Cells.Item(18,2).Value = "Efnb3"
$ws.Cells.Item(18,3).Value = "Ephb2"
$ws.Cells.Item(18,4).Value = "FAPs"
$ws.Cells.Item(18,5).Value = 1
$ws.Cells.Item(18,6).Value = [double]"0.3333333333333333"
$ws.Cells.Item(18,7).Value = [double]"0.155836"
$ws.Cells.Item(18,8).Value = [double]"0.467508"
$ws.Cells.Item(18,9).Value = [double]"0.1119783395680848"
$ws.Cells.Item(18,10).Value = [double]"0.1119783395680848"
$ws.Cells.Item(18,11).Value = 3
$ws.Cells.Item(18,12).Value = 1
$ws.Cells.Item(18,13).Value = [double]"8.445752666666667"
$ws.Cells.Item(18,14).Value = [double]"25.337258"
$ws.Cells.Item(18,15).Value = [double]"0.9204452022087118"
$ws.Cells.Item(18,16).Value = [double]"0.920445202208712"
$ws.Cells.Item(18,17).Value = [double]"1.316152312562667"
$ws.Cells.Item(18,18).Value = [double]"11.845370813064"
$ws.Cells.Item(18,19).Value = [double]"0.1030699254067416"
$ws.Cells.Item(18,20).Value = [double]"0.1030699254067416"

# Row 19 (new)
$ws.Cells.Item(19,1).Value = "Resolving-Mac"
$ws.Cells.Item(19,2).Value = "Efnb3"
$ws.Cells.Item(19,3).Value = "Ephb2"
$ws.Cells.Item(19,4).Value = "Inflammatory-Mac"
$ws.Cells.Item(19,5).Value = 1
$ws.Cells.Item(19,6).Value = [double]"0.3333333333333333"
$ws.Cells.Item(19,7).Value = [double]"0.155836"
$ws.Cells.Item(19,8).Value = [double]"0.467508"
$ws.Cells.Item(19,9).Value = [double]"0.1119783395680848"
$ws.Cells.Item(19,10).Value = [double]"0.1119783395680848"
$ws.Cells.Item(19,11).Value = 2
$ws.Cells.Item(19,12).Value = [double]"0.6666666666666666"
$ws.Cells.Item(19,13).Value = [double]"0.05870933333333334"
$ws.Cells.Item(19,14).Value = [double]"0.176128"
$ws.Cells.Item(19,15).Value = [double]"0.006398331365399365"
$ws.Cells.Item(19,16).Value = [double]"0.006398331365399365"
$ws.Cells.Item(19,17).Value = [double]"0.009149027669333334"
$ws.Cells.Item(19,18).Value = [double]"0.082341249024"
$ws.Cells.Item(19,19).Value = [double]"0.0007164745223038176"
$ws.Cells.Item(19,20).Value = [double]"0.0007164745223038175"

# Row 20 (new)
$ws.Cells.Item(20,1).Value = "Resolving-Mac"
$ws.Cells.Item(20,2).Value = "Efnb3"
$ws.Cells.Item(20,3).Value = "Ephb2"
$ws.Cells.Item(20,4).Value = "MuSCs"
$ws.Cells.Item(20,5).Value = 1
$ws.Cells.Item(20,6).Value = [double]"0.3333333333333333"
$ws.Cells.Item(20,7).Value = [double]"0.155836"
$ws.Cells.Item(20,8).Value = [double]"0.467508"
$ws.Cells.Item(20,9).Value = [double]"0.1119783395680848"
$ws.Cells.Item(20,10).Value = [double]"0.1119783395680848"
$ws.Cells.Item(20,11).Value = 3
$ws.Cells.Item(20,12).Value = 1
$ws.Cells.Item(20,13).Value = [double]"0.4594193333333333"
$ws.Cells.Item(20,14).Value = [double]"1.378258"
$ws.Cells.Item(20,15).Value = [double]"0.05006899181852175"
$ws.Cells.Item(20,16).Value = [double]"0.05006899181852175"
$ws.Cells.Item(20,17).Value = [double]"0.07159407122933334"
$ws.Cells.Item(20,18).Value = [double]"0.644346641064"
$ws.Cells.Item(20,19).Value = [double]"0.005606642567686087"
$ws.Cells.Item(20,20).Value = [double]"0.005606642567686086"

# Row 21 (new)
$ws.Cells.Item(21,1).Value = "Resolving-Mac"
$ws.Cells.Item(21,2).Value = "Efnb3"
$ws.Cells.Item(21,3).Value = "Ephb2"
$ws.Cells.Item(21,4).Value = "Resolving-Mac"
$ws.Cells.Item(21,5).Value = 1
$ws.Cells.Item(21,6).Value = [double]"0.3333333333333333"
$ws.Cells.Item(21,7).Value = [double]"0.155836"
$ws.Cells.Item(21,8).Value = [double]"0.467508"
$ws.Cells.Item(21,9).Value = [double]"0.1119783395680848"
$ws.Cells.Item(21,10).Value = [double]"0.1119783395680848"
$ws.Cells.Item(21,11).Value = 1
$ws.Cells.Item(21,12).Value = [double]"0.3333333333333333"
$ws.Cells.Item(21,13).Value = [double]"0.003065"
$ws.Cells.Item(21,14).Value = [double]"0.009195"
$ws.Cells.Item(21,15).Value = [double]"0.0003340335262130221"
$ws.Cells.Item(21,16).Value = [double]"0.0003340335262130222"
$ws.Cells.Item(21,17).Value = [double]"0.00047763734"
$ws.Cells.Item(21,18).Value = [double]"0.00429873606"
$ws.Cells.Item(21,19).Value = [double]"3.740451962540654e-05"
$ws.Cells.Item(21,20).Value = [double]"3.740451962540654e-05"
